$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.19322657585144
$ws.Range("B1").Value = 2.281124591827393
$ws.Range("C1").Value = 6.526005744934082
$ws.Range("D1").Value = 2.302768707275391
$ws.Range("E1").Value = 1.186935305595398
